{"js": "// 1) The \"city of the future\" paragraph currently ends with a stray\n//    image-markdown fragment: \" [The illustration is by Howard V. Brown.]](images/10k_years.jpg)\".\n//    Strip that fragment off the end of the paragraph, then add a new\n//    paragraph right after it carrying the figure placeholder.\nconst trailer1 = \" [The illustration is by Howard V. Brown.]](images/10k_years.jpg)\";\nconst hits1 = context.document.body.search(trailer1, { matchCase: true });\nhits1.load(\"items\");\nawait context.sync();\n\nif (hits1.items.length > 0) {\n  const hit = hits1.items[0];\n  const paras = hit.paragraphs;\n  paras.load(\"items\");\n  await context.sync();\n  const paragraph = paras.items[0];\n\n  // Remove the trailing fragment from the paragraph text.\n  hit.delete();\n  await context.sync();\n\n  // Insert the new figure-placeholder paragraph right after it.\n  paragraph.insertParagraph(\"[INSERT FIGURE 40.1 NEAR HERE]\", Word.InsertLocation.after);\n  await context.sync();\n}\n\n// 2) The paragraph that is just a leftover image-markdown closer\n//    \"](images/si_192202_cover.jpg)\" becomes the figure placeholder text.\nconst trailer2 = \"](images/si_192202_cover.jpg)\";\nconst hits2 = context.document.body.search(trailer2, { matchCase: true });\nhits2.load(\"items\");\nawait context.sync();\n\nif (hits2.items.length > 0) {\n  const hit2 = hits2.items[0];\n  const paras2 = hit2.paragraphs;\n  paras2.load(\"items\");\n  await context.sync();\n  const paragraph2 = paras2.items[0];\n\n  paragraph2.insertText(\"[INSERT FIGURE 40.2 NEAR HERE]\", Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# 1) The \"city of the future\" paragraph currently ends with a stray\n#    image-markdown fragment: \" [The illustration is by Howard V. Brown.]](images/10k_years.jpg)\".\n#    Strip that fragment off the end of the paragraph, then add a new\n#    paragraph right after it carrying the figure placeholder.\n$marker1 = \" [The illustration is by Howard V. Brown.]](images/10k_years.jpg)\"\n$rng1 = $d.Content\n$found1 = $rng1.Find.Execute($marker1)\nif ($found1) {\n    $paragraph = $rng1.Paragraphs.Item(1)\n    $rng1.Delete()\n\n    $endRng = $paragraph.Range\n    $endRng.Collapse(0)\n    $endRng.InsertParagraphAfter()\n\n    $newParagraph = $d.Paragraphs.Item($paragraph.Index + 1)\n    $newParagraph.Range.Text = \"[INSERT FIGURE 40.1 NEAR HERE]\"\n}\n\n# 2) The paragraph that is just a leftover image-markdown closer\n#    \"](images/si_192202_cover.jpg)\" becomes the figure placeholder text.\n$marker2 = \"](images/si_192202_cover.jpg)\"\n$rng2 = $d.Content\n$found2 = $rng2.Find.Execute($marker2)\nif ($found2) {\n    $rng2.Text = \"[INSERT FIGURE 40.2 NEAR HERE]\"\n}\n"}
